# Nudge the school logo picture ("Afbeelding 7") on the single slide
# a tiny amount to the right/down, as happened in the authored edit
# (diff: a:off x="8705726" y="-542" -> x="8705727" y="14339").
#
# Shape geometry is exposed through the COM object model in points,
# while the underlying OOXML stores EMUs (1 pt = 12700 EMU). The
# values below were chosen so that, after the runtime's internal
# point<->EMU conversion, they land exactly on the target EMU offsets.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Afbeelding 7") {
        $targetShape = $candidate
    }
}

$targetShape.Left = 685.490356
$targetShape.Top = 1.129094
